# "fix error in intro" - day1/Introduction.pptx
#
# Slide 7, "Content Placeholder 2" shape has a paragraph whose text is the
# hyperlinked workshop URL followed by a tab character. The URL text (and
# the part of it after the domain) is corrected from the old
# "2017-September-Microbial-Community-Analysis-Workshop" workshop link text
# to the new "2017_2018-single-cell-RNA-sequencing-Workshop-UCD_UCB_UCSF"
# text, while the underlying hyperlink (rId2) is preserved/reapplied on the
# whole visible URL text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The URL + trailing tab live in the 2nd paragraph of this placeholder.
$para = $tr.Paragraphs(2, 1)

$domainPart = "https://ucdavis-bioinformatics-training.github.io"
$oldPathPart = "/2017-September-Microbial-Community-Analysis-Workshop/"
$newPathPart = "/2017_2018-single-cell-RNA-sequencing-Workshop-UCD_UCB_UCSF/"

# Sanity check - the paragraph should start with the (unchanged) domain,
# immediately followed by the old path that needs replacing.
if ($para.Text.IndexOf($domainPart + $oldPathPart) -eq 0) {
    # Select just the trailing path portion of the URL (leave the domain
    # portion, and the trailing tab run, untouched) and retype it with the
    # corrected workshop path. The hyperlink formatting (a:hlinkClick to
    # rId2) carries over onto the replacement text automatically.
    $pathRange = $para.Characters($domainPart.Length + 1, $oldPathPart.Length)
    $pathRange.Text = $newPathPart
}
